$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column values are stored as plain text (dot-separated thousands,
# e.g. "68.819.23"), so force text number format before assigning numeric-looking
# values, then reset the style back to Normal so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.819.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.734.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.732.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000243"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.361.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.738.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.814.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "496.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.724"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  -8.41%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.880.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.43%  "
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.667.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.325"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "433.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.741.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.83%  "
